$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.564914445839804
$ws.Range("D2").Value = 0.008319282504338332
$ws.Range("E2").Value = 0.06000286272903566
$ws.Range("F2").Value = 4.392711194258879
$ws.Range("G2").Value = 0.002651238078790563
$ws.Range("J2").Value = 0.152571998438205
$ws.Range("K2").Value = 1.72295569776216
$ws.Range("L2").Value = 0.386611203594029
$ws.Range("M2").Value = 0.4095522688283282
$ws.Range("N2").Value = 4.530819572623386
$ws.Range("B3").Value = 1.548440089871121
$ws.Range("D3").Value = 0.007600696708852439
$ws.Range("E3").Value = 0.05923413267756139
$ws.Range("F3").Value = 4.379598386480239
$ws.Range("G3").Value = 0.002656111474630321
$ws.Range("J3").Value = 0.1517105256063491
$ws.Range("K3").Value = 1.664842800434087
$ws.Range("L3").Value = 0.3809137346762554
$ws.Range("M3").Value = 0.4046397912732402
$ws.Range("N3").Value = 4.541032088007753
$ws.Range("B4").Value = 1.53907788036193
$ws.Range("D4").Value = 0.007156417427712114
$ws.Range("E4").Value = 0.05875207406124883
$ws.Range("F4").Value = 4.37316651644592
$ws.Range("G4").Value = 0.002659263458754335
$ws.Range("J4").Value = 0.151166165108739
$ws.Range("K4").Value = 1.630226399633244
$ws.Range("L4").Value = 0.3776276929537943
$ws.Range("M4").Value = 0.4018310531147371
$ws.Range("N4").Value = 4.548183155748617
$ws.Range("B5").Value = 1.53545235972453
$ws.Range("D5").Value = 0.006974559918219825
$ws.Range("E5").Value = 0.05855307515452868
$ws.Range("F5").Value = 4.370952632885732
$ws.Range("G5").Value = 0.002660588206757994
$ws.Range("J5").Value = 0.1509404316411427
$ws.Range("K5").Value = 1.61638739109398
$ws.Range("L5").Value = 0.3763419686890188
$ws.Range("M5").Value = 0.400738706646635
$ws.Range("N5").Value = 4.551318616924391
$ws.Range("B6").Value = 1.534861807875018
$ws.Range("D6").Value = 0.006944312541495634
$ws.Range("E6").Value = 0.05851987634971856
$ws.Range("F6").Value = 4.370609608688625
$ws.Range("G6").Value = 0.002660810617293914
$ws.Range("J6").Value = 0.150902712317599
$ws.Range("K6").Value = 1.614105577004125
$ws.Range("L6").Value = 0.3761316982355822
$ws.Range("M6").Value = 0.4005604798176847
$ws.Range("N6").Value = 4.551852624280343
$ws.Range("B7").Value = 1.539028217042556
$ws.Range("D7").Value = 0.007153968165241054
$ws.Range("E7").Value = 0.05874940067770495
$ws.Range("F7").Value = 4.373135010786811
$ws.Range("G7").Value = 0.00265928116152225
$ws.Range("J7").Value = 0.1511631366218911
$ws.Range("K7").Value = 1.630038679511728
$ws.Range("L7").Value = 0.377610137157987
$ws.Range("M7").Value = 0.4018161097771475
$ws.Range("N7").Value = 4.5482245455906
$ws.Range("B8").Value = 1.559077936913098
$ws.Range("D8").Value = 0.008072126207562036
$ws.Range("E8").Value = 0.05973987476740739
$ws.Range("F8").Value = 4.387853726001993
$ws.Range("G8").Value = 0.002652885359051949
$ws.Range("J8").Value = 0.1522781454050364
$ws.Range("K8").Value = 1.702697040098599
$ws.Range("L8").Value = 0.38460266399386
$ws.Range("M8").Value = 0.4078153972466652
$ws.Range("N8").Value = 4.53415804264985
$ws.Range("B9").Value = 1.604362790230624
$ws.Range("D9").Value = 0.009850166011588612
$ws.Range("E9").Value = 0.06160371216576976
$ws.Range("F9").Value = 4.42957584661994
$ws.Range("G9").Value = 0.002641604375645915
$ws.Range("J9").Value = 0.1543436303248544
$ws.Range("K9").Value = 1.85365902508093
$ws.Range("L9").Value = 0.4000003457611001
$ws.Range("M9").Value = 0.4212257951245846
$ws.Range("N9").Value = 4.513564916995463
$ws.Range("B10").Value = 1.641266535248263
$ws.Range("D10").Value = 0.01114546359183066
$ws.Range("E10").Value = 0.06292712044968418
$ws.Range("F10").Value = 4.46808973606619
$ws.Range("G10").Value = 0.002634076705403627
$ws.Range("J10").Value = 0.1557891565254153
$ws.Range("K10").Value = 1.969793011781917
$ws.Range("L10").Value = 0.4123443900731303
$ws.Range("M10").Value = 0.4320822499454451
$ws.Range("N10").Value = 4.502705069640044
$ws.Range("B11").Value = 1.658843307082861
$ws.Range("D11").Value = 0.0117329195806164
$ws.Range("E11").Value = 0.06351960861490191
$ws.Range("F11").Value = 4.487322977250216
$ws.Range("G11").Value = 0.002630815517098477
$ws.Range("J11").Value = 0.1564315266177161
$ws.Range("K11").Value = 2.023771863950856
$ws.Range("L11").Value = 0.4181849711717831
$ws.Range("M11").Value = 0.4372393576488278
$ws.Range("N11").Value = 4.498693720684841
$ws.Range("B12").Value = 1.665612481122764
$ws.Range("D12").Value = 0.01195516382917106
$ws.Range("E12").Value = 0.06374262885520388
$ws.Range("F12").Value = 4.494852740992656
$ws.Range("G12").Value = 0.002629603920526305
$ws.Range("J12").Value = 0.1566726186927028
$ws.Range("K12").Value = 2.044378196132357
$ws.Range("L12").Value = 0.4204290799491588
$ws.Range("M12").Value = 0.4392236214589715
$ws.Range("N12").Value = 4.497308443719589
$ws.Range("B13").Value = 1.664149586047813
$ws.Range("D13").Value = 0.01190730830985132
$ws.Range("E13").Value = 0.06369465668190166
$ws.Range("F13").Value = 4.493220102610138
$ws.Range("G13").Value = 0.002629863823219682
$ws.Range("J13").Value = 0.1566207908122976
$ws.Range("K13").Value = 2.039932875599391
$ws.Range("L13").Value = 0.4199443292641405
$ws.Range("M13").Value = 0.4387948796765215
$ws.Range("N13").Value = 4.497600837339363
$ws.Range("B14").Value = 1.659397943097957
$ws.Range("D14").Value = 0.0117512077582802
$ws.Range("E14").Value = 0.06353798335506156
$ws.Range("F14").Value = 4.487937512895797
$ws.Range("G14").Value = 0.002630715371053015
$ws.Range("J14").Value = 0.1564514045319285
$ws.Range("K14").Value = 2.025463835269022
$ws.Range("L14").Value = 0.4183689457227615
$ws.Range("M14").Value = 0.4374019755327936
$ws.Range("N14").Value = 4.498577071920067
$ws.Range("B15").Value = 1.656502163837928
$ws.Range("D15").Value = 0.011655565280698
$ws.Range("E15").Value = 0.06344184250000495
$ws.Range("F15").Value = 4.484733889151101
$ws.Range("G15").Value = 0.002631240005613769
$ws.Range("J15").Value = 0.1563473701193843
$ws.Range("K15").Value = 2.016622725980085
$ws.Range("L15").Value = 0.4174081989195173
$ws.Range("M15").Value = 0.4365528670222929
$ws.Range("N15").Value = 4.499192464904155
$ws.Range("B16").Value = 1.640133709763603
$ws.Range("D16").Value = 0.01110703982990913
$ws.Range("E16").Value = 0.0628882105411126
$ws.Range("F16").Value = 4.466867288088537
$ws.Range("G16").Value = 0.00263429310480539
$ws.Range("J16").Value = 0.1557468724480664
$ws.Range("K16").Value = 1.966288528386258
$ws.Range("L16").Value = 0.4119672280644835
$ws.Range("M16").Value = 0.4317496136872165
$ws.Range("N16").Value = 4.5029859256598
$ws.Range("B17").Value = 1.630294124441349
$ws.Range("D17").Value = 0.01077011182884036
$ws.Range("E17").Value = 0.06254615336116487
$ws.Range("F17").Value = 4.456345616370925
$ws.Range("G17").Value = 0.002636207791665851
$ws.Range("J17").Value = 0.1553746110113234
$ws.Range("K17").Value = 1.935704744953227
$ws.Range("L17").Value = 0.4086870661925417
$ws.Range("M17").Value = 0.428858899323231
$ws.Range("N17").Value = 4.505551108364145
$ws.Range("B18").Value = 1.62470894027706
$ws.Range("D18").Value = 0.01057614944686591
$ws.Range("E18").Value = 0.06234851169863909
$ws.Range("F18").Value = 4.450455058793125
$ws.Range("G18").Value = 0.00263732443600602
$ws.Range("J18").Value = 0.155159063085673
$ws.Range("K18").Value = 1.918221872646086
$ws.Range("L18").Value = 0.4068216005361904
$ws.Range("M18").Value = 0.4272168005420127
$ws.Range("N18").Value = 4.507113951129099
$ws.Range("B19").Value = 1.622830658683625
$ws.Range("D19").Value = 0.01051044654650113
$ws.Range("E19").Value = 0.06228143826493948
$ws.Range("F19").Value = 4.4484883022171
$ws.Range("G19").Value = 0.002637705155408952
$ws.Range("J19").Value = 0.155085835317065
$ws.Range("K19").Value = 1.912321025663687
$ws.Range("L19").Value = 0.4061936258348595
$ws.Range("M19").Value = 0.4266643470043405
$ws.Range("N19").Value = 4.507658111881184
$ws.Range("B20").Value = 1.631333878612509
$ws.Range("D20").Value = 0.01080599570085639
$ws.Range("E20").Value = 0.06258265880614644
$ws.Range("F20").Value = 4.457448978167093
$ws.Range("G20").Value = 0.002636002380421015
$ws.Range("J20").Value = 0.1554143869881717
$ws.Range("K20").Value = 1.938949248273218
$ws.Range("L20").Value = 0.4090340510968673
$ws.Range("M20").Value = 0.4291644927968292
$ws.Range("N20").Value = 4.505268991769086
$ws.Range("B21").Value = 1.660790544581062
$ws.Range("D21").Value = 0.01179706369340749
$ws.Range("E21").Value = 0.06358403830108728
$ws.Range("F21").Value = 4.489482444412573
$ws.Range("G21").Value = 0.002630464618432007
$ws.Range("J21").Value = 0.1565012157804091
$ws.Range("K21").Value = 2.029709241714102
$ws.Range("L21").Value = 0.4188307943401526
$ws.Range("M21").Value = 0.43781025388882
$ws.Range("N21").Value = 4.49828669730816
$ws.Range("B22").Value = 1.680702075967815
$ws.Range("D22").Value = 0.01244356539174873
$ws.Range("E22").Value = 0.06423068822669364
$ws.Range("F22").Value = 4.511855371969517
$ws.Range("G22").Value = 0.002626981383061884
$ws.Range("J22").Value = 0.1571989465171981
$ws.Range("K22").Value = 2.089992284371533
$ws.Range("L22").Value = 0.4254224398420519
$ws.Range("M22").Value = 0.4436436451028882
$ws.Range("N22").Value = 4.494502948346536
$ws.Range("B23").Value = 1.67001461273145
$ws.Range("D23").Value = 0.0120986122369473
$ws.Range("E23").Value = 0.06388626384191021
$ws.Range("F23").Value = 4.499782941312986
$ws.Range("G23").Value = 0.002628828046181442
$ws.Range("J23").Value = 0.1568276960486763
$ws.Range("K23").Value = 2.057729519423049
$ws.Range("L23").Value = 0.4218870616362977
$ws.Range("M23").Value = 0.4405135310935648
$ws.Range("N23").Value = 4.496451018010532
$ws.Range("B24").Value = 1.630863582320444
$ws.Range("D24").Value = 0.01078977340972642
$ws.Range("E24").Value = 0.06256615777616759
$ws.Range("F24").Value = 4.456949654621127
$ws.Range("G24").Value = 0.002636095197315075
$ws.Range("J24").Value = 0.1553964090373903
$ws.Range("K24").Value = 1.937482096401482
$ws.Range("L24").Value = 0.4088771158571376
$ws.Range("M24").Value = 0.4290262722750313
$ws.Range("N24").Value = 4.505396262302298
$ws.Range("B25").Value = 1.591473714900786
$ws.Range("D25").Value = 0.00937129021493277
$ws.Range("E25").Value = 0.0611077019046764
$ws.Range("F25").Value = 4.416910097319615
$ws.Range("G25").Value = 0.002644522032556852
$ws.Range("J25").Value = 0.1537976342835314
$ws.Range("K25").Value = 1.811906617872239
$ws.Range("L25").Value = 0.3956540380251994
$ws.Range("M25").Value = 0.4174216858865947
$ws.Range("N25").Value = 4.518386648135746
